# Rename worksheets (drop the "StreamData" prefix, use underscore before period number)
$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("StreamDataPeriod1").Name = "Period_1"
$wb.Worksheets.Item("StreamDataPeriod2").Name = "Period_2"
$wb.Worksheets.Item("StreamDataPeriod3").Name = "Period_3"

# Update the selected cell on the (3rd) sheet that is currently tab-selected
# (was StreamDataPeriod3, now Period_3): move selection from E13 to C19.
$ws3 = $wb.Worksheets.Item("Period_3")
$ws3.Activate()
$ws3.Range("C19").Select()
